$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "2014-10"
$ws.Cells.Item(2, 2).Value = 100.463
$ws.Cells.Item(2, 3).Value = 100.8533
$ws.Cells.Item(2, 4).Value = 98.51179999999999
$ws.Cells.Item(2, 5).Value = 97.14709999999999
$ws.Cells.Item(2, 6).Value = 97.5778
$ws.Cells.Item(2, 7).Value = 100.1513
$ws.Cells.Item(2, 8).Value = 96.73869999999999
$ws.Cells.Item(2, 9).Value = 101.6072

$ws.Cells.Item(3, 1).Value = "2014-11"
$ws.Cells.Item(3, 2).Value = 100.418
$ws.Cells.Item(3, 3).Value = 101.2526
$ws.Cells.Item(3, 4).Value = 98.715
$ws.Cells.Item(3, 5).Value = 97.2957
$ws.Cells.Item(3, 6).Value = 97.10550000000001
$ws.Cells.Item(3, 7).Value = 99.4469
$ws.Cells.Item(3, 8).Value = 96.98690000000001
$ws.Cells.Item(3, 9).Value = 101.7578

$ws.Cells.Item(4, 1).Value = "2014-12"
$ws.Cells.Item(4, 2).Value = 100.9294
$ws.Cells.Item(4, 3).Value = 101.3678
$ws.Cells.Item(4, 4).Value = 98.76819999999999
$ws.Cells.Item(4, 5).Value = 97.8236
$ws.Cells.Item(4, 6).Value = 97.70359999999999
$ws.Cells.Item(4, 7).Value = 99.6678
$ws.Cells.Item(4, 8).Value = 97.6891
$ws.Cells.Item(4, 9).Value = 102.1332

$ws.Cells.Item(5, 1).Value = "2014-01"
$ws.Cells.Item(5, 2).Value = 99.78360000000001
$ws.Cells.Item(5, 3).Value = 98.2961
$ws.Cells.Item(5, 4).Value = 97.214
$ws.Cells.Item(5, 5).Value = 95.04770000000001
$ws.Cells.Item(5, 6).Value = 95.2283
$ws.Cells.Item(5, 7).Value = 99.3723
$ws.Cells.Item(5, 8).Value = 96.39579999999999
$ws.Cells.Item(5, 9).Value = 98.9318

$ws.Cells.Item(6, 1).Value = "2014-02"
$ws.Cells.Item(6, 2).Value = 98.4059
$ws.Cells.Item(6, 3).Value = 98.4761
$ws.Cells.Item(6, 4).Value = 97.152
$ws.Cells.Item(6, 5).Value = 95.31229999999999
$ws.Cells.Item(6, 6).Value = 96.0157
$ws.Cells.Item(6, 7).Value = 99.9436
$ws.Cells.Item(6, 8).Value = 96.5812
$ws.Cells.Item(6, 9).Value = 99.083

$ws.Cells.Item(7, 1).Value = "2014-03"
$ws.Cells.Item(7, 2).Value = 98.51819999999999
$ws.Cells.Item(7, 3).Value = 98.2959
$ws.Cells.Item(7, 4).Value = 97.381
$ws.Cells.Item(7, 5).Value = 96.0236
$ws.Cells.Item(7, 6).Value = 96.5509
$ws.Cells.Item(7, 7).Value = 100.2399
$ws.Cells.Item(7, 8).Value = 96.9764
$ws.Cells.Item(7, 9).Value = 99.0642

$ws.Cells.Item(8, 1).Value = "2014-04"
$ws.Cells.Item(8, 2).Value = 99.9986
$ws.Cells.Item(8, 3).Value = 98.4555
$ws.Cells.Item(8, 4).Value = 97.90470000000001
$ws.Cells.Item(8, 5).Value = 95.85129999999999
$ws.Cells.Item(8, 6).Value = 96.82429999999999
$ws.Cells.Item(8, 7).Value = 99.6738
$ws.Cells.Item(8, 8).Value = 96.6865
$ws.Cells.Item(8, 9).Value = 99.3507

$ws.Cells.Item(9, 1).Value = "2014-05"
$ws.Cells.Item(9, 2).Value = 99.7745
$ws.Cells.Item(9, 3).Value = 98.90389999999999
$ws.Cells.Item(9, 4).Value = 98.31059999999999
$ws.Cells.Item(9, 5).Value = 96.4669
$ws.Cells.Item(9, 6).Value = 97.5936
$ws.Cells.Item(9, 7).Value = 99.7897
$ws.Cells.Item(9, 8).Value = 96.92919999999999
$ws.Cells.Item(9, 9).Value = 99.7176

$ws.Cells.Item(10, 1).Value = "2014-06"
$ws.Cells.Item(10, 2).Value = 99.43810000000001
$ws.Cells.Item(10, 3).Value = 99.4105
$ws.Cells.Item(10, 4).Value = 98.44580000000001
$ws.Cells.Item(10, 5).Value = 97.3593
$ws.Cells.Item(10, 6).Value = 96.9041
$ws.Cells.Item(10, 7).Value = 99.874
$ws.Cells.Item(10, 8).Value = 96.8279
$ws.Cells.Item(10, 9).Value = 99.5898

$ws.Cells.Item(11, 1).Value = "2014-07"
$ws.Cells.Item(11, 2).Value = 100.0073
$ws.Cells.Item(11, 3).Value = 99.61020000000001
$ws.Cells.Item(11, 4).Value = 98.47880000000001
$ws.Cells.Item(11, 5).Value = 96.9389
$ws.Cells.Item(11, 6).Value = 96.6087
$ws.Cells.Item(11, 7).Value = 100.2555
$ws.Cells.Item(11, 8).Value = 96.95489999999999
$ws.Cells.Item(11, 9).Value = 100.6999

$ws.Cells.Item(12, 1).Value = "2014-08"
$ws.Cells.Item(12, 2).Value = 99.6854
$ws.Cells.Item(12, 3).Value = 99.5655
$ws.Cells.Item(12, 4).Value = 98.36150000000001
$ws.Cells.Item(12, 5).Value = 97.0741
$ws.Cells.Item(12, 6).Value = 97.52760000000001
$ws.Cells.Item(12, 7).Value = 100.4844
$ws.Cells.Item(12, 8).Value = 96.9573
$ws.Cells.Item(12, 9).Value = 101.4498

$ws.Cells.Item(13, 1).Value = "2014-09"
$ws.Cells.Item(13, 2).Value = 100.1714
$ws.Cells.Item(13, 3).Value = 100.0744
$ws.Cells.Item(13, 4).Value = 98.4264
$ws.Cells.Item(13, 5).Value = 97.4152
$ws.Cells.Item(13, 6).Value = 97.3015
$ws.Cells.Item(13, 7).Value = 100.398
$ws.Cells.Item(13, 8).Value = 96.648
$ws.Cells.Item(13, 9).Value = 101.2954

$ws.Cells.Item(14, 1).Value = "2015-10"
$ws.Cells.Item(14, 2).Value = 101.2
$ws.Cells.Item(14, 3).Value = 99
$ws.Cells.Item(14, 4).Value = 99.8
$ws.Cells.Item(14, 5).Value = 96.2
$ws.Cells.Item(14, 6).Value = 98.09999999999999
$ws.Cells.Item(14, 7).Value = 98.7
$ws.Cells.Item(14, 8).Value = 98.2
$ws.Cells.Item(14, 9).Value = 99.8

$ws.Cells.Item(15, 1).Value = "2015-11"
$ws.Cells.Item(15, 2).Value = 100.9847
$ws.Cells.Item(15, 3).Value = 98.8154
$ws.Cells.Item(15, 4).Value = 99.8091
$ws.Cells.Item(15, 5).Value = 96.548
$ws.Cells.Item(15, 6).Value = 97.60299999999999
$ws.Cells.Item(15, 7).Value = 98.5506
$ws.Cells.Item(15, 8).Value = 98.0133
$ws.Cells.Item(15, 9).Value = 99.3124

$ws.Cells.Item(16, 1).Value = "2015-12"
$ws.Cells.Item(16, 2).Value = 100.7381
$ws.Cells.Item(16, 3).Value = 98.92359999999999
$ws.Cells.Item(16, 4).Value = 100.0548
$ws.Cells.Item(16, 5).Value = 96.23990000000001
$ws.Cells.Item(16, 6).Value = 96.8963
$ws.Cells.Item(16, 7).Value = 98.2748
$ws.Cells.Item(16, 8).Value = 97.3741
$ws.Cells.Item(16, 9).Value = 99.0222

$ws.Cells.Item(17, 1).Value = "2015-01"
$ws.Cells.Item(17, 2).Value = 99.8904
$ws.Cells.Item(17, 3).Value = 100.7932
$ws.Cells.Item(17, 4).Value = 98.85380000000001
$ws.Cells.Item(17, 5).Value = 97.1773
$ws.Cells.Item(17, 6).Value = 98.4014
$ws.Cells.Item(17, 7).Value = 98.9855
$ws.Cells.Item(17, 8).Value = 97.6981
$ws.Cells.Item(17, 9).Value = 102.0026

$ws.Cells.Item(18, 1).Value = "2015-02"
$ws.Cells.Item(18, 2).Value = 101.6073
$ws.Cells.Item(18, 3).Value = 100.133
$ws.Cells.Item(18, 4).Value = 99.0026
$ws.Cells.Item(18, 5).Value = 97.35509999999999
$ws.Cells.Item(18, 6).Value = 98.36879999999999
$ws.Cells.Item(18, 7).Value = 98.67529999999999
$ws.Cells.Item(18, 8).Value = 97.7517
$ws.Cells.Item(18, 9).Value = 102.1189

$ws.Cells.Item(19, 1).Value = "2015-03"
$ws.Cells.Item(19, 2).Value = 101.6609
$ws.Cells.Item(19, 3).Value = 100.2589
$ws.Cells.Item(19, 4).Value = 99.04730000000001
$ws.Cells.Item(19, 5).Value = 96.9932
$ws.Cells.Item(19, 6).Value = 98.3302
$ws.Cells.Item(19, 7).Value = 98.934
$ws.Cells.Item(19, 8).Value = 97.8502
$ws.Cells.Item(19, 9).Value = 102.2948

$ws.Cells.Item(20, 1).Value = "2015-04"
$ws.Cells.Item(20, 2).Value = 100.9775
$ws.Cells.Item(20, 3).Value = 100.131
$ws.Cells.Item(20, 4).Value = 99.2512
$ws.Cells.Item(20, 5).Value = 97.16379999999999
$ws.Cells.Item(20, 6).Value = 98.4089
$ws.Cells.Item(20, 7).Value = 99.10129999999999
$ws.Cells.Item(20, 8).Value = 98.1292
$ws.Cells.Item(20, 9).Value = 102.2021

$ws.Cells.Item(21, 1).Value = "2015-05"
$ws.Cells.Item(21, 2).Value = 101.3414
$ws.Cells.Item(21, 3).Value = 100.4651
$ws.Cells.Item(21, 4).Value = 99.2497
$ws.Cells.Item(21, 5).Value = 96.9944
$ws.Cells.Item(21, 6).Value = 97.9498
$ws.Cells.Item(21, 7).Value = 98.8639
$ws.Cells.Item(21, 8).Value = 98.00069999999999
$ws.Cells.Item(21, 9).Value = 101.8562

$ws.Cells.Item(22, 1).Value = "2015-06"
$ws.Cells.Item(22, 2).Value = 101.5088
$ws.Cells.Item(22, 3).Value = 100.083
$ws.Cells.Item(22, 4).Value = 99.0989
$ws.Cells.Item(22, 5).Value = 96.6802
$ws.Cells.Item(22, 6).Value = 98.79049999999999
$ws.Cells.Item(22, 7).Value = 99.2
$ws.Cells.Item(22, 8).Value = 97.96210000000001
$ws.Cells.Item(22, 9).Value = 101.0808

$ws.Cells.Item(23, 1).Value = "2015-07"
$ws.Cells.Item(23, 2).Value = 100.7357
$ws.Cells.Item(23, 3).Value = 100.1063
$ws.Cells.Item(23, 4).Value = 99.1382
$ws.Cells.Item(23, 5).Value = 96.7242
$ws.Cells.Item(23, 6).Value = 98.7761
$ws.Cells.Item(23, 7).Value = 98.6738
$ws.Cells.Item(23, 8).Value = 97.92140000000001
$ws.Cells.Item(23, 9).Value = 101.1151

$ws.Cells.Item(24, 1).Value = "2015-08"
$ws.Cells.Item(24, 2).Value = 101.5334
$ws.Cells.Item(24, 3).Value = 99.3245
$ws.Cells.Item(24, 4).Value = 99.07729999999999
$ws.Cells.Item(24, 5).Value = 96.5262
$ws.Cells.Item(24, 6).Value = 97.9483
$ws.Cells.Item(24, 7).Value = 98.02160000000001
$ws.Cells.Item(24, 8).Value = 97.8378
$ws.Cells.Item(24, 9).Value = 100.2926

$ws.Cells.Item(25, 1).Value = "2015-09"
$ws.Cells.Item(25, 2).Value = 101.4548
$ws.Cells.Item(25, 3).Value = 99.3377
$ws.Cells.Item(25, 4).Value = 99.5977
$ws.Cells.Item(25, 5).Value = 96.9545
$ws.Cells.Item(25, 6).Value = 98.5474
$ws.Cells.Item(25, 7).Value = 98.6442
$ws.Cells.Item(25, 8).Value = 98.414
$ws.Cells.Item(25, 9).Value = 100.0327

$ws.Cells.Item(26, 1).Value = "2016-10"
$ws.Cells.Item(26, 2).Value = 98.90000000000001
$ws.Cells.Item(26, 3).Value = 100.4
$ws.Cells.Item(26, 4).Value = 99.3
$ws.Cells.Item(26, 5).Value = 99.7
$ws.Cells.Item(26, 6).Value = 97.5
$ws.Cells.Item(26, 7).Value = 99.3
$ws.Cells.Item(26, 8).Value = 97.40000000000001
$ws.Cells.Item(26, 9).Value = 100.7

$ws.Cells.Item(27, 1).Value = "2016-11"
$ws.Cells.Item(27, 2).Value = 99.2
$ws.Cells.Item(27, 3).Value = 100.9
$ws.Cells.Item(27, 4).Value = 99.90000000000001
$ws.Cells.Item(27, 5).Value = 100.7
$ws.Cells.Item(27, 6).Value = 98.40000000000001
$ws.Cells.Item(27, 7).Value = 100
$ws.Cells.Item(27, 8).Value = 97
$ws.Cells.Item(27, 9).Value = 97

$ws.Cells.Item(28, 1).Value = "2016-12"
$ws.Cells.Item(28, 2).Value = 99
$ws.Cells.Item(28, 3).Value = 101.6
$ws.Cells.Item(28, 4).Value = 100
$ws.Cells.Item(28, 5).Value = 101.5
$ws.Cells.Item(28, 6).Value = 99.5
$ws.Cells.Item(28, 7).Value = 100.7
$ws.Cells.Item(28, 8).Value = 97.40000000000001
$ws.Cells.Item(28, 9).Value = 96.59999999999999

$ws.Cells.Item(29, 1).Value = "2016-01"
$ws.Cells.Item(29, 2).Value = 102.0032
$ws.Cells.Item(29, 3).Value = 99.5857
$ws.Cells.Item(29, 4).Value = 99.5638
$ws.Cells.Item(29, 5).Value = 96.62009999999999
$ws.Cells.Item(29, 6).Value = 95.6564
$ws.Cells.Item(29, 7).Value = 99.30629999999999
$ws.Cells.Item(29, 8).Value = 97.86620000000001
$ws.Cells.Item(29, 9).Value = 100.0436

$ws.Cells.Item(30, 1).Value = "2016-02"
$ws.Cells.Item(30, 2).Value = 101.7057
$ws.Cells.Item(30, 3).Value = 99.82680000000001
$ws.Cells.Item(30, 4).Value = 99.7364
$ws.Cells.Item(30, 5).Value = 96.4996
$ws.Cells.Item(30, 6).Value = 94.9855
$ws.Cells.Item(30, 7).Value = 99.479
$ws.Cells.Item(30, 8).Value = 97.9098
$ws.Cells.Item(30, 9).Value = 99.5421

$ws.Cells.Item(31, 1).Value = "2016-03"
$ws.Cells.Item(31, 2).Value = 101.4516
$ws.Cells.Item(31, 3).Value = 99.9516
$ws.Cells.Item(31, 4).Value = 99.7236
$ws.Cells.Item(31, 5).Value = 96.5654
$ws.Cells.Item(31, 6).Value = 94.70820000000001
$ws.Cells.Item(31, 7).Value = 99.2667
$ws.Cells.Item(31, 8).Value = 97.81229999999999
$ws.Cells.Item(31, 9).Value = 99.0185

$ws.Cells.Item(32, 1).Value = "2016-04"
$ws.Cells.Item(32, 2).Value = 101.325
$ws.Cells.Item(32, 3).Value = 99.93729999999999
$ws.Cells.Item(32, 4).Value = 99.6767
$ws.Cells.Item(32, 5).Value = 96.41
$ws.Cells.Item(32, 6).Value = 94.30800000000001
$ws.Cells.Item(32, 7).Value = 99.0039
$ws.Cells.Item(32, 8).Value = 97.84569999999999
$ws.Cells.Item(32, 9).Value = 99.5338

$ws.Cells.Item(33, 1).Value = "2016-05"
$ws.Cells.Item(33, 2).Value = 101
$ws.Cells.Item(33, 3).Value = 100.1
$ws.Cells.Item(33, 4).Value = 99.2
$ws.Cells.Item(33, 5).Value = 97.7
$ws.Cells.Item(33, 6).Value = 94.59999999999999
$ws.Cells.Item(33, 7).Value = 99.2
$ws.Cells.Item(33, 8).Value = 97.59999999999999
$ws.Cells.Item(33, 9).Value = 99.8

$ws.Cells.Item(34, 1).Value = "2016-06"
$ws.Cells.Item(34, 2).Value = 101
$ws.Cells.Item(34, 3).Value = 100.9
$ws.Cells.Item(34, 4).Value = 99.3
$ws.Cells.Item(34, 5).Value = 98.3
$ws.Cells.Item(34, 6).Value = 94.90000000000001
$ws.Cells.Item(34, 7).Value = 99.2
$ws.Cells.Item(34, 8).Value = 97.7
$ws.Cells.Item(34, 9).Value = 100.6

$ws.Cells.Item(35, 1).Value = "2016-07"
$ws.Cells.Item(35, 2).Value = 100.8
$ws.Cells.Item(35, 3).Value = 101.1
$ws.Cells.Item(35, 4).Value = 99.5
$ws.Cells.Item(35, 5).Value = 98.3
$ws.Cells.Item(35, 6).Value = 95.09999999999999
$ws.Cells.Item(35, 7).Value = 99.8
$ws.Cells.Item(35, 8).Value = 97.8
$ws.Cells.Item(35, 9).Value = 101.2

$ws.Cells.Item(36, 1).Value = "2016-08"
$ws.Cells.Item(36, 2).Value = 100.5
$ws.Cells.Item(36, 3).Value = 101.2
$ws.Cells.Item(36, 4).Value = 99.5
$ws.Cells.Item(36, 5).Value = 98.90000000000001
$ws.Cells.Item(36, 6).Value = 96
$ws.Cells.Item(36, 7).Value = 100
$ws.Cells.Item(36, 8).Value = 97.8
$ws.Cells.Item(36, 9).Value = 101

$ws.Cells.Item(37, 1).Value = "2016-09"
$ws.Cells.Item(37, 2).Value = 99.40000000000001
$ws.Cells.Item(37, 3).Value = 100.9
$ws.Cells.Item(37, 4).Value = 99.2
$ws.Cells.Item(37, 5).Value = 98.7
$ws.Cells.Item(37, 6).Value = 96.3
$ws.Cells.Item(37, 7).Value = 99.3
$ws.Cells.Item(37, 8).Value = 97.40000000000001
$ws.Cells.Item(37, 9).Value = 101.5

$ws.Cells.Item(38, 1).Value = "2017-10"
$ws.Cells.Item(38, 2).Value = 99
$ws.Cells.Item(38, 3).Value = 102
$ws.Cells.Item(38, 4).Value = 98.2
$ws.Cells.Item(38, 5).Value = 99.2
$ws.Cells.Item(38, 6).Value = 100.6
$ws.Cells.Item(38, 7).Value = 100.4
$ws.Cells.Item(38, 8).Value = 98.40000000000001
$ws.Cells.Item(38, 9).Value = 96.5

$ws.Cells.Item(39, 1).Value = "2017-11"
$ws.Cells.Item(39, 2).Value = 99.09999999999999
$ws.Cells.Item(39, 3).Value = 101.6
$ws.Cells.Item(39, 4).Value = 97.7
$ws.Cells.Item(39, 5).Value = 98.3
$ws.Cells.Item(39, 6).Value = 100.2
$ws.Cells.Item(39, 7).Value = 100.3
$ws.Cells.Item(39, 8).Value = 98.59999999999999
$ws.Cells.Item(39, 9).Value = 99.2

$ws.Cells.Item(40, 1).Value = "2017-12"
$ws.Cells.Item(40, 2).Value = 99.2
$ws.Cells.Item(40, 3).Value = 101
$ws.Cells.Item(40, 4).Value = 97.2
$ws.Cells.Item(40, 5).Value = 97.59999999999999
$ws.Cells.Item(40, 6).Value = 99.2
$ws.Cells.Item(40, 7).Value = 99.40000000000001
$ws.Cells.Item(40, 8).Value = 98.59999999999999
$ws.Cells.Item(40, 9).Value = 100.7

$ws.Cells.Item(41, 1).Value = "2017-01"
$ws.Cells.Item(41, 2).Value = 99.2
$ws.Cells.Item(41, 3).Value = 102
$ws.Cells.Item(41, 4).Value = 100.2
$ws.Cells.Item(41, 5).Value = 101.9
$ws.Cells.Item(41, 6).Value = 100
$ws.Cells.Item(41, 7).Value = 100.8
$ws.Cells.Item(41, 8).Value = 97.3
$ws.Cells.Item(41, 9).Value = 96.59999999999999

$ws.Cells.Item(42, 1).Value = "2017-02"
$ws.Cells.Item(42, 2).Value = 99.3
$ws.Cells.Item(42, 3).Value = 102.6
$ws.Cells.Item(42, 4).Value = 100.2
$ws.Cells.Item(42, 5).Value = 101.5
$ws.Cells.Item(42, 6).Value = 100.5
$ws.Cells.Item(42, 7).Value = 100.6
$ws.Cells.Item(42, 8).Value = 97.09999999999999
$ws.Cells.Item(42, 9).Value = 97.59999999999999

$ws.Cells.Item(43, 1).Value = "2017-03"
$ws.Cells.Item(43, 2).Value = 99.40000000000001
$ws.Cells.Item(43, 3).Value = 102.8
$ws.Cells.Item(43, 4).Value = 100.4
$ws.Cells.Item(43, 5).Value = 101.9
$ws.Cells.Item(43, 6).Value = 100.9
$ws.Cells.Item(43, 7).Value = 100.9
$ws.Cells.Item(43, 8).Value = 97.59999999999999
$ws.Cells.Item(43, 9).Value = 97.2

$ws.Cells.Item(44, 1).Value = "2017-04"
$ws.Cells.Item(44, 2).Value = 99.3
$ws.Cells.Item(44, 3).Value = 103
$ws.Cells.Item(44, 4).Value = 99.8
$ws.Cells.Item(44, 5).Value = 101.9
$ws.Cells.Item(44, 6).Value = 102.2
$ws.Cells.Item(44, 7).Value = 101.4
$ws.Cells.Item(44, 8).Value = 98.09999999999999
$ws.Cells.Item(44, 9).Value = 96.7

$ws.Cells.Item(45, 1).Value = "2017-05"
$ws.Cells.Item(45, 2).Value = 99.5
$ws.Cells.Item(45, 3).Value = 102.9
$ws.Cells.Item(45, 4).Value = 100.4
$ws.Cells.Item(45, 5).Value = 101.1
$ws.Cells.Item(45, 6).Value = 102.7
$ws.Cells.Item(45, 7).Value = 101.5
$ws.Cells.Item(45, 8).Value = 98.2
$ws.Cells.Item(45, 9).Value = 96.40000000000001

$ws.Cells.Item(46, 1).Value = "2017-06"
$ws.Cells.Item(46, 2).Value = 99.2
$ws.Cells.Item(46, 3).Value = 102.2
$ws.Cells.Item(46, 4).Value = 100.1
$ws.Cells.Item(46, 5).Value = 100.6
$ws.Cells.Item(46, 6).Value = 102.6
$ws.Cells.Item(46, 7).Value = 101.3
$ws.Cells.Item(46, 8).Value = 98.2
$ws.Cells.Item(46, 9).Value = 96.09999999999999

$ws.Cells.Item(47, 1).Value = "2017-07"
$ws.Cells.Item(47, 2).Value = 99.40000000000001
$ws.Cells.Item(47, 3).Value = 101.8
$ws.Cells.Item(47, 4).Value = 99.8
$ws.Cells.Item(47, 5).Value = 100.6
$ws.Cells.Item(47, 6).Value = 102.5
$ws.Cells.Item(47, 7).Value = 100.9
$ws.Cells.Item(47, 8).Value = 98.09999999999999
$ws.Cells.Item(47, 9).Value = 96.09999999999999

$ws.Cells.Item(48, 1).Value = "2017-08"
$ws.Cells.Item(48, 2).Value = 98.90000000000001
$ws.Cells.Item(48, 3).Value = 101.7
$ws.Cells.Item(48, 4).Value = 99.40000000000001
$ws.Cells.Item(48, 5).Value = 99.8
$ws.Cells.Item(48, 6).Value = 101.9
$ws.Cells.Item(48, 7).Value = 100.8
$ws.Cells.Item(48, 8).Value = 98.3
$ws.Cells.Item(48, 9).Value = 96.2

$ws.Cells.Item(49, 1).Value = "2017-09"
$ws.Cells.Item(49, 2).Value = 98.90000000000001
$ws.Cells.Item(49, 3).Value = 101.6
$ws.Cells.Item(49, 4).Value = 98.2
$ws.Cells.Item(49, 5).Value = 99.40000000000001
$ws.Cells.Item(49, 6).Value = 101.2
$ws.Cells.Item(49, 7).Value = 100.1
$ws.Cells.Item(49, 8).Value = 98.5
$ws.Cells.Item(49, 9).Value = 95.59999999999999
